# Fuel Prod Imp Exp Balancing Priorities.xlsx — apply the authored edit:
#  - About!C1 "last updated" date bumped from 1/3/2024 to 3/28/2024
#  - FPIEBP!B3:D3 (hard coal) priority order changed from 3,2,1 to 1,3,2
#  - Selection on FPIEBP moved from F4 to E3 (matches saved view state)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")

# Update the "last updated" date on the About sheet.
$wsAbout.Range("C1").Value = "3/28/2024"

# Update the hard coal priority ranking (production / imports / exports).
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# Restore the saved cursor position on the FPIEBP sheet.
$wsFpiebp.Activate()
[void]$wsFpiebp.Range("E3").Select()
